# Append a new data row (row 22) to the NIFTY_Options_Analysis sheet,
# mirroring the structure/formatting of the last existing row (row 21),
# then fill in the new day's values (2026-01-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NIFTY_Options_Analysis")

$srcRow = 21
$dstRow = 22

# Copy the source row's cell range (values + number formats/fills/fonts/borders)
# into the new row first, so the new row looks exactly like the rows above it.
$ws.Range("A" + $srcRow + ":AE" + $srcRow).Copy($ws.Range("A" + $dstRow + ":AE" + $dstRow))

# A (Date) and E (Position_Size) look like a date / a percentage to Excel's
# automatic value parser, so force them to Text first to keep them as the
# literal strings "2026-01-20" and "100%" (matching how the source data is
# stored as plain text in this sheet).
$ws.Cells.Item($dstRow, 1).NumberFormat = "@"
$ws.Cells.Item($dstRow, 1).Value = "2026-01-20"    # A - Date

$ws.Cells.Item($dstRow, 2).Value = "10:00:07"      # B - Time
$ws.Cells.Item($dstRow, 3).Value = "AVOID"         # C - Signal
$ws.Cells.Item($dstRow, 4).Value = "AVOID"         # D - Signal_Tier

$ws.Cells.Item($dstRow, 5).NumberFormat = "@"
$ws.Cells.Item($dstRow, 5).Value = "100%"          # E - Position_Size

$ws.Cells.Item($dstRow, 6).Value = "TRADEABLE"     # F - Premium_Quality
$ws.Cells.Item($dstRow, 7).Value = 0                # G - Total_Score
$ws.Cells.Item($dstRow, 8).Value = 25458            # H - NIFTY_Spot
$ws.Cells.Item($dstRow, 9).Value = 12.3             # I - VIX
$ws.Cells.Item($dstRow, 10).Value = 0.98            # J - VIX_Trend
$ws.Cells.Item($dstRow, 11).Value = 0               # K - VIX_Score
$ws.Cells.Item($dstRow, 12).Value = 45.9            # L - IV_Rank
$ws.Cells.Item($dstRow, 13).Value = "UNKNOWN"      # M - Market_Regime
$ws.Cells.Item($dstRow, 14).Value = 0               # N - Regime_Score
$ws.Cells.Item($dstRow, 15).Value = "UNKNOWN"      # O - OI_Pattern
$ws.Cells.Item($dstRow, 16).Value = 0               # P - OI_Score
$ws.Cells.Item($dstRow, 17).Value = 0               # Q - Theta_Score
$ws.Cells.Item($dstRow, 18).Value = 0               # R - Gamma_Score
$ws.Cells.Item($dstRow, 19).Value = 0               # S - Vega_Score
$ws.Cells.Item($dstRow, 20).Value = "NONE"         # T - Best_Strategy
$ws.Cells.Item($dstRow, 21).Value = ""              # U - Expiry_1 (blank)
$ws.Cells.Item($dstRow, 22).Value = 0               # V - Days_To_Expiry_1
$ws.Cells.Item($dstRow, 23).Value = 0               # W - Straddle_Premium
$ws.Cells.Item($dstRow, 24).Value = 0               # X - Straddle_Theta
$ws.Cells.Item($dstRow, 25).Value = 0               # Y - Straddle_Gamma
$ws.Cells.Item($dstRow, 26).Value = 0               # Z - Strangle_Premium
$ws.Cells.Item($dstRow, 27).Value = 0               # AA - Strangle_Theta
$ws.Cells.Item($dstRow, 28).Value = 0               # AB - Strangle_Gamma
$ws.Cells.Item($dstRow, 29).Value = "HARD VETO: CPR TRENDING DAY: Price 25458.00 below BC 25573.82 - BEARISH TRENDING DAY likely"  # AC - Recommendation
$ws.Cells.Item($dstRow, 30).Value = "CPR TRENDING DAY: Price 25458.00 below BC 25573.82 - BEARISH TRENDING DAY likely"             # AD - Risk_Factors
$ws.Cells.Item($dstRow, 31).Value = "Yes"          # AE - Telegram_Sent
